$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 77.127561
$ws.Range("H2").Value = 231.382683
$ws.Range("I2").Value = 0.2899056040435161
$ws.Range("J2").Value = 0.2899056040435161
$ws.Range("M2").Value = 4.959409333333333
$ws.Range("N2").Value = 14.878228
$ws.Range("O2").Value = 0.8271666313262851
$ws.Range("P2").Value = 0.8271666313262852
$ws.Range("Q2").Value = 382.507145880636
$ws.Range("R2").Value = 3442.564312925724
$ws.Range("S2").Value = 0.2398002418992871
$ws.Range("T2").Value = 0.2398002418992871
$ws.Range("G3").Value = 77.127561
$ws.Range("H3").Value = 231.382683
$ws.Range("I3").Value = 0.2899056040435161
$ws.Range("J3").Value = 0.2899056040435161
$ws.Range("O3").Value = 0.09421438109281059
$ws.Range("P3").Value = 0.09421438109281059
$ws.Range("Q3").Value = 43.567610984184
$ws.Range("R3").Value = 392.108498857656
$ws.Range("S3").Value = 0.02731327706029728
$ws.Range("T3").Value = 0.02731327706029728
$ws.Range("G4").Value = 77.127561
$ws.Range("H4").Value = 231.382683
$ws.Range("I4").Value = 0.2899056040435161
$ws.Range("J4").Value = 0.2899056040435161
$ws.Range("O4").Value = 0.07861898758090437
$ws.Range("P4").Value = 0.07861898758090438
$ws.Range("Q4").Value = 36.355824102066
$ws.Range("R4").Value = 327.202416918594
$ws.Range("S4").Value = 0.02279208508393177
$ws.Range("T4").Value = 0.02279208508393178
$ws.Range("I5").Value = 0.443028781054351
$ws.Range("J5").Value = 0.443028781054351
$ws.Range("M5").Value = 4.959409333333333
$ws.Range("N5").Value = 14.878228
$ws.Range("O5").Value = 0.8271666313262851
$ws.Range("P5").Value = 0.8271666313262852
$ws.Range("Q5").Value = 584.5408719958378
$ws.Range("R5").Value = 5260.86784796254
$ws.Range("S5").Value = 0.3664586244053178
$ws.Range("T5").Value = 0.3664586244053178
$ws.Range("I6").Value = 0.443028781054351
$ws.Range("J6").Value = 0.443028781054351
$ws.Range("O6").Value = 0.09421438109281059
$ws.Range("P6").Value = 0.09421438109281059
$ws.Range("S6").Value = 0.04173968241333797
$ws.Range("T6").Value = 0.04173968241333797
$ws.Range("I7").Value = 0.443028781054351
$ws.Range("J7").Value = 0.443028781054351
$ws.Range("O7").Value = 0.07861898758090437
$ws.Range("P7").Value = 0.07861898758090438
$ws.Range("S7").Value = 0.03483047423569522
$ws.Range("T7").Value = 0.03483047423569522
$ws.Range("I8").Value = 0.267065614902133
$ws.Range("J8").Value = 0.2670656149021329
$ws.Range("M8").Value = 4.959409333333333
$ws.Range("N8").Value = 14.878228
$ws.Range("O8").Value = 0.8271666313262851
$ws.Range("P8").Value = 0.8271666313262852
$ws.Range("Q8").Value = 352.3716157751062
$ws.Range("R8").Value = 3171.344541975956
$ws.Range("S8").Value = 0.2209077650216803
$ws.Range("T8").Value = 0.2209077650216802
$ws.Range("I9").Value = 0.267065614902133
$ws.Range("J9").Value = 0.2670656149021329
$ws.Range("O9").Value = 0.09421438109281059
$ws.Range("P9").Value = 0.09421438109281059
$ws.Range("S9").Value = 0.02516142161917535
$ws.Range("T9").Value = 0.02516142161917535
$ws.Range("I10").Value = 0.267065614902133
$ws.Range("J10").Value = 0.2670656149021329
$ws.Range("O10").Value = 0.07861898758090437
$ws.Range("P10").Value = 0.07861898758090438
$ws.Range("S10").Value = 0.02099642826127738
$ws.Range("T10").Value = 0.02099642826127738
